# Update the violations list: the original three fines (FORD MUSTANG,
# BMW 330I, NISSAN SENTRA) are replaced by a single new fine entry for a
# CADILLAC ESCALADE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A2")
$text = $cell.Value2

$oldBlock = "FORD MUSTANG, 2024, Black`nDate and Time of Issuing The Fine`n12 Jul 2025, 4:58 am`nAmount`nAED 700`nSource`nDubai Police`nBlack points`n-`nBMW 330I 2023, Black`nDate and Time of Issuing The Fine`n12 Jul 2025, 3:41 am`nAmount`nAED 600`nSource`nDubai Police`nBlack points`n-`nNISSAN SENTRA, 2023, Blue`nDate and Time of Issuing The Fine`n11 Jul 2025, 1:18 pm`nAmount`nAED 400`nSource`nDubai Police`nBlack points"

$newBlock = "CADILLAC ESCALADE, 2023, Blue`nDate and Time of Issuing The Fine`n14 Jul 2025, 12:10 am`nAmount`nAED 600`nSource`nDubai Police`nBlack points"

$text = $text.Replace($oldBlock, $newBlock)

$cell.Value = $text
